# Update the cryptos list worksheet with refreshed prices / volumes,
# and restore the original Bitcoin / ShibaInu / InternetComputer row order swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (these columns hold inline strings like "42.979.75" or "  -0.97%  ",
    # some of which are also syntactically valid numbers and would otherwise
    # be auto-converted). Resetting the style back to Normal afterwards
    # drops the transient quote-prefix formatting so the cell keeps using
    # the sheet default style, matching the original workbook.
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell 'D2' '42.979.75'
Set-TextCell 'E2' '  -0.97%  '
Set-TextCell 'D3' '2.337.68'
Set-TextCell 'E3' '  +1.16%  '
Set-TextCell 'E4' '  +0.07%  '
Set-TextCell 'D5' '307.06'
Set-TextCell 'D6' '101.51'
Set-TextCell 'E6' '  -2.04%  '
Set-TextCell 'D7' '0.509'
Set-TextCell 'E7' '  -4.40%  '
Set-TextCell 'E8' '  +0.04%  '
Set-TextCell 'D9' '0.510'
Set-TextCell 'E9' '  -3.74%  '
Set-TextCell 'D10' '34.84'
Set-TextCell 'E10' '  -4.87%  '
Set-TextCell 'D11' '52.50'
Set-TextCell 'E11' '  +1.34%  '
Set-TextCell 'D12' '0.0797'
Set-TextCell 'E12' '  -2.29%  '
Set-TextCell 'E13' '  +0.74%  '
Set-TextCell 'D14' '6.85'
Set-TextCell 'E14' '  -2.91%  '
Set-TextCell 'D15' '15.80'
Set-TextCell 'E15' '  +4.63%  '
Set-TextCell 'D16' '2.341.82'
Set-TextCell 'E16' '  -4.28%  '
Set-TextCell 'D17' '0.827'
Set-TextCell 'E17' '  +2.06%  '
Set-TextCell 'D18' '42.905.42'
Set-TextCell 'E18' '  -0.92%  '
Set-TextCell 'B19' 'InternetComputer(DFINITY)'
Set-TextCell 'C19' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D19' '11.74'
Set-TextCell 'E19' '  -5.06%  '
Set-TextCell 'B20' 'ShibaInu'
Set-TextCell 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D20' '0.0₃0911'
Set-TextCell 'E20' '  -2.40%  '
Set-TextCell 'D21' '6.19'
Set-TextCell 'E21' '  -0.15%  '
Set-TextCell 'D22' '68.51'
Set-TextCell 'E22' '  +0.51%  '
Set-TextCell 'D23' '236.75'
Set-TextCell 'E23' '  -2.24%  '
Set-TextCell 'E24' '  -0.79%  '
Set-TextCell 'D25' '2.56'
Set-TextCell 'E25' '  -2.56%  '
Set-TextCell 'E26' '  -0.06%  '
Set-TextCell 'D27' '25.62'
Set-TextCell 'E27' '  +3.05%  '
Set-TextCell 'E28' '  +1.08%  '
Set-TextCell 'D29' '35.27'
Set-TextCell 'E29' '  -4.38%  '
Set-TextCell 'D30' '9.30'
Set-TextCell 'E30' '  -3.87%  '
Set-TextCell 'D31' '161.69'
Set-TextCell 'E31' '  -3.94%  '
Set-TextCell 'E32' '  +0.05%  '
Set-TextCell 'E33' '  -3.45%  '
Set-TextCell 'D34' '4.65'
Set-TextCell 'E34' '  +5.13%  '
Set-TextCell 'E35' '  -3.31%  '
Set-TextCell 'D36' '2.45'
Set-TextCell 'E36' '  -3.59%  '
Set-TextCell 'D37' '0.0725'
Set-TextCell 'E37' '  -2.62%  '
Set-TextCell 'E38' '  -1.84%  '
Set-TextCell 'D39' '2.91'
Set-TextCell 'E39' '  -4.99%  '
Set-TextCell 'E40' '  -3.89%  '
Set-TextCell 'E41' '  -2.74%  '
Set-TextCell 'E42' '  +5.66%  '
Set-TextCell 'D43' '2.023.23'
Set-TextCell 'E43' '  +2.01%  '
Set-TextCell 'E44' '  -3.94%  '
Set-TextCell 'D45' '18.90'
Set-TextCell 'E45' '  -1.93%  '
Set-TextCell 'D46' '10.19'
Set-TextCell 'E46' '  +2.03%  '
Set-TextCell 'D47' '2.94'
Set-TextCell 'E47' '  -3.06%  '
Set-TextCell 'D48' '55.88'
Set-TextCell 'E48' '  -0.01%  '
Set-TextCell 'D50' '2.561.91'
Set-TextCell 'E50' '  +1.02%  '
Set-TextCell 'E51' '  +0.95%  '
